$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Host $ws.Name
Write-Host $wb.Worksheets.Count
